$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UseCase3")

# Step 1: update columns I (minimal) and J (condition) on the existing 12 rows
$ws.Cells.Item(53, 9).Value = "minimal"
$ws.Cells.Item(53, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(54, 9).Value = "minimal"
$ws.Cells.Item(54, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(55, 9).Value = "minimal"
$ws.Cells.Item(55, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(56, 9).Value = "minimal"
$ws.Cells.Item(56, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(57, 9).Value = "minimal"
$ws.Cells.Item(57, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(58, 9).Value = "minimal"
$ws.Cells.Item(58, 10).Value = "([histopat] = '0')"
$ws.Cells.Item(59, 9).Value = "minimal"
$ws.Cells.Item(59, 10).Value = "([histopat] = '0' AND [tnm_pt] = '')"
$ws.Cells.Item(60, 9).Value = "minimal"
$ws.Cells.Item(60, 10).Value = "([histopat] = '0' AND [tnm_pt] = '')"
$ws.Cells.Item(61, 9).Value = "minimal"
$ws.Cells.Item(61, 10).Value = "([histopat] = '0' AND [tnm_pn] = '')"
$ws.Cells.Item(62, 9).Value = "minimal"
$ws.Cells.Item(62, 10).Value = "([histopat] = '0' AND [tnm_pn] = '')"
$ws.Cells.Item(63, 9).Value = "minimal"
$ws.Cells.Item(63, 10).Value = "([histopat] = '0' AND [tnm_pm] = '')"
$ws.Cells.Item(64, 9).Value = "minimal"
$ws.Cells.Item(64, 10).Value = "([histopat] = '0' AND [tnm_pm] = '')"

# Step 2: insert a new "minimal_req" row after each "permissible" row, top to bottom
# Block pt: type row 53, permissible row 54, new row 55
$ws.Rows.Item(53).Copy()
$ws.Rows.Item(55).Insert()
$ws.Cells.Item(55, 1).Borders.LineStyle = 1
$ws.Cells.Item(55, 4).Value = "uc3_tnm_pt_min_req"
$ws.Cells.Item(55, 5).Value = "completeness"
$ws.Cells.Item(55, 6).Value = "minimal_req"
$ws.Cells.Item(55, 7).Value = "[histopat] = '0' AND [tnm_pt] = ''"
$ws.Cells.Item(55, 14).Value = "pT classification of the primary tumour is expected to meet minimal requirements"

# Block pn: type row 56, permissible row 57, new row 58
$ws.Rows.Item(56).Copy()
$ws.Rows.Item(58).Insert()
$ws.Cells.Item(58, 1).Borders.LineStyle = 1
$ws.Cells.Item(58, 4).Value = "uc3_tnm_pn_min_req"
$ws.Cells.Item(58, 5).Value = "completeness"
$ws.Cells.Item(58, 6).Value = "minimal_req"
$ws.Cells.Item(58, 7).Value = "[histopat] = '0' AND [tnm_pn] = ''"
$ws.Cells.Item(58, 14).Value = "pN classification of the primary tumour is expected to meet minimal requirements"

# Block pm: type row 59, permissible row 60, new row 61
$ws.Rows.Item(59).Copy()
$ws.Rows.Item(61).Insert()
$ws.Cells.Item(61, 1).Borders.LineStyle = 1
$ws.Cells.Item(61, 4).Value = "uc3_tnm_pm_min_req"
$ws.Cells.Item(61, 5).Value = "completeness"
$ws.Cells.Item(61, 6).Value = "minimal_req"
$ws.Cells.Item(61, 7).Value = "[histopat] = '0' AND [tnm_pm] = ''"
$ws.Cells.Item(61, 14).Value = "pM classification of the primary tumour is expected to meet minimal requirements"

# Block ct: type row 62, permissible row 63, new row 64
$ws.Rows.Item(62).Copy()
$ws.Rows.Item(64).Insert()
$ws.Cells.Item(64, 1).Borders.LineStyle = 1
$ws.Cells.Item(64, 4).Value = "uc3_tnm_ct_min_req"
$ws.Cells.Item(64, 5).Value = "completeness"
$ws.Cells.Item(64, 6).Value = "minimal_req"
$ws.Cells.Item(64, 7).Value = "[histopat] = '0' AND [tnm_pt] = '' AND [tnm_ct] = ''"
$ws.Cells.Item(64, 14).Value = "cT classification of the primary tumour is expected to meet minimal requirements"

# Block cn: type row 65, permissible row 66, new row 67
$ws.Rows.Item(65).Copy()
$ws.Rows.Item(67).Insert()
$ws.Cells.Item(67, 1).Borders.LineStyle = 1
$ws.Cells.Item(67, 4).Value = "uc3_tnm_cn_min_req"
$ws.Cells.Item(67, 5).Value = "completeness"
$ws.Cells.Item(67, 6).Value = "minimal_req"
$ws.Cells.Item(67, 7).Value = "[histopat] = '0' AND [tnm_pn] = '' AND [tnm_cn] = ''"
$ws.Cells.Item(67, 14).Value = "cN classification of the primary tumour is expected to meet minimal requirements"

# Block cm: type row 68, permissible row 69, new row 70
$ws.Rows.Item(68).Copy()
$ws.Rows.Item(70).Insert()
$ws.Cells.Item(70, 1).Borders.LineStyle = 1
$ws.Cells.Item(70, 4).Value = "uc3_tnm_cm_min_req"
$ws.Cells.Item(70, 5).Value = "completeness"
$ws.Cells.Item(70, 6).Value = "minimal_req"
$ws.Cells.Item(70, 7).Value = "[histopat] = '0' AND [tnm_pm] = '' AND [tnm_cm] = ''"
$ws.Cells.Item(70, 14).Value = "cM classification of the primary tumour is expected to meet minimal requirements"

# Step 2b: re-number column A (sequential id = row-2) for rows 53-70
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(63, 1).Value = 61
$ws.Cells.Item(64, 1).Value = 62
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(66, 1).Value = 64
$ws.Cells.Item(67, 1).Value = 65
$ws.Cells.Item(68, 1).Value = 66
$ws.Cells.Item(69, 1).Value = 67
$ws.Cells.Item(70, 1).Value = 68

# Step 3: fix typo her2ihq -> her2ihc on UseCase7 sheet
$ws7 = $wb.Worksheets.Item("UseCase7")
$ws7.Cells.Item(64, 10).Value = "([her2ihc] = '2')"
$ws7.Cells.Item(65, 10).Value = "([her2ihc] = '2')"
$ws7.Cells.Item(66, 7).Value = "[her2ihc] = '2' AND [her2fish] = ''"
$ws7.Cells.Item(66, 10).Value = "([her2ihc] = '2')"

